# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (well outside the A1:E51 data range) used to force text-typed
# values into Price/Volume columns without Excel re-interpreting them as numbers.
$scratch = $ws.Range("ZZ1")

function Set-TextValue($cellRef, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "70.602.36"
Set-TextValue "E2" "  -0.12%  "

Set-TextValue "D3" "3.513.20"
Set-TextValue "E3" "  -1.95%  "

Set-TextValue "E4" "  -0.13%  "

Set-TextValue "D5" "623.38"
Set-TextValue "E5" "  +4.20%  "

Set-TextValue "D6" "172.57"
Set-TextValue "E6" "  -0.31%  "

Set-TextValue "D7" "0.609"
Set-TextValue "E7" "  -1.31%  "

Set-TextValue "D8" "3.507.35"
Set-TextValue "E8" "  -1.98%  "

Set-TextValue "E9" "  +0.02%  "

Set-TextValue "D10" "0.198"
Set-TextValue "E10" "  -0.54%  "

Set-TextValue "D11" "7.14"
Set-TextValue "E11" "  -3.78%  "

Set-TextValue "D12" "0.587"
Set-TextValue "E12" "  -0.54%  "

Set-TextValue "D13" "46.35"
Set-TextValue "E13" "  -0.91%  "

Set-TextValue "E14" "  -0.86%  "

Set-TextValue "D15" "4.087.60"
Set-TextValue "E15" "  -1.80%  "

Set-TextValue "E16" "  -0.42%  "

Set-TextValue "D17" "609.38"

Set-TextValue "D18" "3.521.92"
Set-TextValue "E18" "  -1.78%  "

Set-TextValue "D19" "70.701.79"
Set-TextValue "E19" "  -0.09%  "

Set-TextValue "E20" "  +1.23%  "

Set-TextValue "D21" "17.75"
Set-TextValue "E21" "  +1.54%  "

Set-TextValue "D22" "0.881"
Set-TextValue "E22" "  -0.66%  "

Set-TextValue "D23" "9.09"
Set-TextValue "E23" "  -2.53%  "

Set-TextValue "E24" "  -2.05%  "

Set-TextValue "D25" "97.31"
Set-TextValue "E25" "  +0.22%  "

Set-TextValue "E26" "  -1.02%  "

Set-TextValue "E27" "  -0.02%  "

Set-TextValue "E28" "  -2.95%  "

Set-TextValue "D29" "33.51"
Set-TextValue "E29" "  -1.16%  "

Set-TextValue "D30" "9.06"
Set-TextValue "E30" "  -1.33%  "

Set-TextValue "D31" "3.01"
Set-TextValue "E31" "  -1.60%  "

Set-TextValue "E32" "  -4.09%  "

Set-TextValue "E33" "  -0.79%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D34" "6.82"
Set-TextValue "E34" "  -5.05%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D35" "627.02"
Set-TextValue "E35" "  -3.23%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D36" "0.0493"
Set-TextValue "E36" "  +2.46%  "

$ws.Range("B37").Value = "Cosmos"
$ws.Range("C37").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D37" "10.81"
Set-TextValue "E37" "  -0.41%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.0993"
Set-TextValue "E38" "  -1.94%  "

Set-TextValue "E39" "  -7.21%  "

Set-TextValue "D40" "56.65"
Set-TextValue "E40" "  -1.16%  "

Set-TextValue "E41" "  +0.19%  "

Set-TextValue "D42" "0.142"
Set-TextValue "E42" "  +0.21%  "

Set-TextValue "D43" "3.343.85"
Set-TextValue "E43" "  -1.51%  "

Set-TextValue "E44" "  +1.20%  "

Set-TextValue "E45" "  +0.21%  "

Set-TextValue "D46" "0.311"
Set-TextValue "E46" "  -3.93%  "

Set-TextValue "D47" "31.99"
Set-TextValue "E47" "  -2.94%  "

Set-TextValue "D48" "2.52"
Set-TextValue "E48" "  -5.15%  "

Set-TextValue "E49" "  -0.58%  "

Set-TextValue "D50" "132.92"
Set-TextValue "E50" "  +0.12%  "

Set-TextValue "D51" "0.155"
Set-TextValue "E51" "  +5.34%  "

$scratch.Clear()
